$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric-looking "Price" (column D) values, forcing text storage
# to match workbook convention (values like "65.754.25" are text, not numbers).
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D42", "D45", "D48", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "65.754.25"
$ws.Range("D3").Value = "3.406.42"
$ws.Range("D5").Value = "562.83"
$ws.Range("D6").Value = "176.03"
$ws.Range("D7").Value = "0.630"
$ws.Range("D8").Value = "3.396.74"
$ws.Range("D11").Value = "0.634"
$ws.Range("D12").Value = "55.10"
$ws.Range("D14").Value = "9.18"
$ws.Range("D15").Value = "3.957.24"
$ws.Range("D16").Value = "18.36"
$ws.Range("D17").Value = "3.405.13"
$ws.Range("D18").Value = "0.119"
$ws.Range("D19").Value = "65.603.29"
$ws.Range("D20").Value = "11.89"
$ws.Range("D21").Value = "0.993"
$ws.Range("D22").Value = "471.58"
$ws.Range("D23").Value = "5.25"
$ws.Range("D25").Value = "86.64"
$ws.Range("D26").Value = "13.54"
$ws.Range("D27").Value = "10.91"
$ws.Range("D29").Value = "8.92"
$ws.Range("D30").Value = "31.05"
$ws.Range("D31").Value = "6.71"
$ws.Range("D32").Value = "11.58"
$ws.Range("D33").Value = "62.75"
$ws.Range("D34").Value = "578.14"
$ws.Range("D42").Value = "3.096.59"
$ws.Range("D45").Value = "0.0416"
$ws.Range("D48").Value = "3.21"
$ws.Range("D50").Value = "8.38"
$ws.Range("D51").Value = "136.40"

foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}

# Update Coin name / Link (columns B, C) and Volume(1h) percentage (column E).
$ws.Range("E2").Value = "  +3.35%  "
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("E5").Value = "  +2.52%  "
$ws.Range("E6").Value = "  +2.25%  "
$ws.Range("E7").Value = "  +2.98%  "
$ws.Range("E8").Value = "  +2.01%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("E10").Value = "  +13.86%  "
$ws.Range("E11").Value = "  +3.51%  "
$ws.Range("E12").Value = "  +3.15%  "
$ws.Range("E13").Value = "  +5.87%  "
$ws.Range("E14").Value = "  +2.93%  "
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("E16").Value = "  +2.97%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("E17").Value = "  +1.83%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("E18").Value = "  +1.84%  "
$ws.Range("E19").Value = "  +3.07%  "
$ws.Range("E20").Value = "  +1.76%  "
$ws.Range("E21").Value = "  +2.10%  "
$ws.Range("E22").Value = "  +13.64%  "
$ws.Range("E23").Value = "  +21.56%  "
$ws.Range("E24").Value = "  +2.73%  "
$ws.Range("E25").Value = "  +4.13%  "
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("E27").Value = "  +3.34%  "
$ws.Range("E28").Value = "  +6.41%  "
$ws.Range("E29").Value = "  +4.11%  "
$ws.Range("E30").Value = "  +6.73%  "
$ws.Range("E31").Value = "  +4.83%  "
$ws.Range("E32").Value = "  +2.23%  "
$ws.Range("E33").Value = "  +9.00%  "
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("E35").Value = "  +2.59%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  -3.20%  "
$ws.Range("E38").Value = "  +5.21%  "
$ws.Range("E39").Value = "  +1.89%  "
$ws.Range("E40").Value = "  +3.01%  "
$ws.Range("E41").Value = "  +2.27%  "
$ws.Range("E42").Value = "  -1.61%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("E44").Value = "  +0.99%  "
$ws.Range("E45").Value = "  +3.54%  "
$ws.Range("E46").Value = "  +3.81%  "
$ws.Range("E47").Value = "  +6.75%  "
$ws.Range("E48").Value = "  -0.89%  "
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("E50").Value = "  +5.15%  "
$ws.Range("E51").Value = "  +1.65%  "
